# Regenerate sval data to filter save games.
# Updates columns B (TB), C (d2S), D (K), E (IP) with new computed values,
# and recomputes G (sum) as B+C+D+E for each data row (rows 2-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    3 = @(1.445647641019636, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    4 = @(3.272327238179451, 1.626987699542094, 0.1496068669990043, 0.5333859586016987)
    5 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987)
    6 = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 2797.565817734744)
    7 = @(0.1169995834814548, 0.3048912486333797, 0.1496068669990043, 0.5333859586016987)
    8 = @(3.272327238179451, 1.626987699542094, 0.7210945179870265, 0.5333859586016987)
    9 = @(0.6545652718822623, 0.3048912486333797, 3.223369029078222, 0.5333859586016987)
}

foreach ($row in $data.Keys | Sort-Object) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = ($b + $c + $d + $e)
}
